# Update the dSF column (column F) values for the fried_max 2023 save-data sheet.
# These figures were repulled/recalculated ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 5
    4  = 5
    5  = 1
    6  = 5
    7  = 2
    8  = -1
    9  = 7
    11 = 7
    12 = -3
    13 = 4
    14 = 2
    16 = 3
    17 = -1
    18 = 5
    19 = 3
    20 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
